# Fix the "Recorded By" (column G) entries so the system-generated
# "System" marker is listed after the human recorder's email instead of
# before it, e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count
$changed = 0

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $current = $cell.Value()
    if ($current -eq $oldText) {
        $cell.Value = $newText
        $changed++
    }
}

Write-Host "Updated $changed cell(s) in column G."
